$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'28.423.60"
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.87%  '

$c = $ws.Range("D3")
$c.Value = "'1.864.35"
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.40%  '

$c = $ws.Range("D4")
$c.Value = "'1.007"
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

$c = $ws.Range("D5")
$c.Value = "'324.98"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("E6").Value = '  +0.16%  '

$c = $ws.Range("D7")
$c.Value = "'0.4560"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.72%  '

$c = $ws.Range("D8")
$c.Value = "'0.3837"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.76%  '

$c = $ws.Range("D9")
$c.Value = "'0.07830"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.38%  '

$c = $ws.Range("D10")
$c.Value = "'0.9879"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.61%  '

$ws.Range("E11").Value = '  -2.50%  '

$c = $ws.Range("D12")
$c.Value = "'1.952.52"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +4.70%  '

$c = $ws.Range("D13")
$c.Value = "'6.901"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.42%  '

$c = $ws.Range("D14")
$c.Value = "'5.637"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.99%  '

$c = $ws.Range("D15")
$c.Value = "'0.06959"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.90%  '

$c = $ws.Range("D16")
$c.Value = "'86.66"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.10%  '

$ws.Range("E17").Value = '  +0.19%  '

$c = $ws.Range("D18")
$c.Value = "'0.000009961"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.02%  '

$c = $ws.Range("D19")
$c.Value = "'16.64"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.53%  '

$c = $ws.Range("D20")
$c.Value = "'1.006"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '

$c = $ws.Range("D21")
$c.Value = "'28.445.13"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.84%  '

$c = $ws.Range("D22")
$c.Value = "'5.250"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.97%  '

$ws.Range("E23").Value = '  -1.17%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D24")
$c.Value = "'2.102"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c = $ws.Range("D25")
$c.Value = "'2.101.18"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.13%  '

$c = $ws.Range("D26")
$c.Value = "'153.36"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.83%  '

$c = $ws.Range("D27")
$c.Value = "'19.12"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.25%  '

$c = $ws.Range("D28")
$c.Value = "'5.641"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.47%  '

$c = $ws.Range("D29")
$c.Value = "'1.942"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.28%  '

$c = $ws.Range("D30")
$c.Value = "'117.59"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.10%  '

$c = $ws.Range("D31")
$c.Value = "'0.09266"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.09%  '

$c = $ws.Range("D32")
$c.Value = "'0.9083"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.31%  '

$c = $ws.Range("D33")
$c.Value = "'5.265"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.39%  '

$ws.Range("E34").Value = '  -0.63%  '

$c = $ws.Range("D35")
$c.Value = "'3.303"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.86%  '

$c = $ws.Range("D36")
$c.Value = "'0.05719"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.67%  '

$c = $ws.Range("D37")
$c.Value = "'1.137"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("E38").Value = '  -3.04%  '

$c = $ws.Range("D39")
$c.Value = "'7.680"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.12%  '

$c = $ws.Range("D40")
$c.Value = "'0.5559"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.60%  '

$ws.Range("E41").Value = '  +0.40%  '

$c = $ws.Range("D42")
$c.Value = "'9.644"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.58%  '

$c = $ws.Range("D43")
$c.Value = "'0.07090"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.47%  '

$c = $ws.Range("D44")
$c.Value = "'11.62"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("E45").Value = '  -1.07%  '

$c = $ws.Range("D46")
$c.Value = "'2.141"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("E47").Value = '  -1.50%  '

$ws.Range("E48").Value = '  -1.29%  '

$c = $ws.Range("D49")
$c.Value = "'111.45"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.22%  '

$c = $ws.Range("D50")
$c.Value = "'2.406"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +3.40%  '

$ws.Range("E51").Value = '  +0.20%  '

Write-Host "Applied all cell updates"